# Update "想去人数" (wanted-to-go count) figures for two rows that appear
# identically on both the "展览" sheet and the "全部类型" sheet.
#
#   展览 sheet:     F4 106 -> 108 ; F5 2711 -> 2731
#   全部类型 sheet: F4 106 -> 108 ; F5 2711 -> 2731

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 108
    $ws.Range("F5").Value = 2731
}
